# Atualizacao de bases das ligas, do dia: 17-05-2024 as 13:59
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2,3,4: same-date fixtures re-ordered by source scraper (3-way rotation) ---

    # Row 2
    $ws.Range("B2").Value = 6720873
    $ws.Range("E2").Value = 'Sportivo Luqueno'
    $ws.Range("F2").Value = 'Sportivo Trinidense'
    $ws.Range("G2").Value = 2
    $ws.Range("H2").Value = 2
    $ws.Range("I2").Value = 'D'
    $ws.Range("J2").Value = 2.625
    $ws.Range("K2").Value = 3.1
    $ws.Range("L2").Value = 2.5
    $ws.Range("M2").Value = 2.3
    $ws.Range("N2").Value = 3.1
    $ws.Range("O2").Value = 2.9
    $ws.Range("P2").Value = -0.25
    $ws.Range("Q2").Value = 2.025
    $ws.Range("R2").Value = 1.775
    $ws.Range("S2").Value = 2.5
    $ws.Range("T2").Value = 1.95
    $ws.Range("U2").Value = 1.85
    $ws.Range("V2").Value = -1
    $ws.Range("W2").Value = 2.1
    $ws.Range("X2").Value = -1
    $ws.Range("Y2").Value = -0.5
    $ws.Range("Z2").Value = 0.3875
    $ws.Range("AA2").Value = 0.95
    $ws.Range("AB2").Value = -1

    # Row 3
    $ws.Range("B3").Value = 6720844
    $ws.Range("E3").Value = 'Guarani Asuncion'
    $ws.Range("F3").Value = 'Olimpia Asuncion'
    $ws.Range("G3").Value = 1
    $ws.Range("H3").Value = 2
    $ws.Range("I3").Value = 'A'
    $ws.Range("J3").Value = 2.45
    $ws.Range("K3").Value = 3
    $ws.Range("L3").Value = 2.75
    $ws.Range("M3").Value = 4
    $ws.Range("N3").Value = 3.2
    $ws.Range("O3").Value = 1.85
    $ws.Range("P3").Value = 0.5
    $ws.Range("Q3").Value = 1.875
    $ws.Range("R3").Value = 1.925
    $ws.Range("S3").Value = 2.5
    $ws.Range("T3").Value = 1.925
    $ws.Range("U3").Value = 1.875
    $ws.Range("V3").Value = -1
    $ws.Range("W3").Value = -1
    $ws.Range("X3").Value = 0.8500000000000001
    $ws.Range("Y3").Value = -1
    $ws.Range("Z3").Value = 0.925
    $ws.Range("AA3").Value = 0.925
    $ws.Range("AB3").Value = -1

    # Row 4
    $ws.Range("B4").Value = 6720843
    $ws.Range("E4").Value = 'Cerro Porteno'
    $ws.Range("F4").Value = 'Libertad Asuncion'
    $ws.Range("G4").Value = 2
    $ws.Range("H4").Value = 0
    $ws.Range("I4").Value = 'H'
    $ws.Range("J4").Value = 2.375
    $ws.Range("K4").Value = 3.2
    $ws.Range("L4").Value = 2.7
    $ws.Range("M4").Value = 3.75
    $ws.Range("N4").Value = 3.3
    $ws.Range("O4").Value = 1.85
    $ws.Range("P4").Value = 0.5
    $ws.Range("Q4").Value = 1.9
    $ws.Range("R4").Value = 1.9
    $ws.Range("S4").Value = 2.5
    $ws.Range("T4").Value = 1.925
    $ws.Range("U4").Value = 1.875
    $ws.Range("V4").Value = 2.75
    $ws.Range("W4").Value = -1
    $ws.Range("X4").Value = -1
    $ws.Range("Y4").Value = 0.8999999999999999
    $ws.Range("Z4").Value = -1
    $ws.Range("AA4").Value = -1
    $ws.Range("AB4").Value = 0.875


    # --- Rows 143,145: same-date fixtures swapped by source scraper ---

    # Row 143
    $ws.Range("B143").Value = 7493311
    $ws.Range("E143").Value = 'General Caballero JLM'
    $ws.Range("F143").Value = 'Olimpia Asuncion'
    $ws.Range("G143").Value = 0
    $ws.Range("H143").Value = 1
    $ws.Range("I143").Value = 'A'
    $ws.Range("J143").Value = 3.4
    $ws.Range("K143").Value = 3.3
    $ws.Range("L143").Value = 2
    $ws.Range("M143").Value = 3.2
    $ws.Range("N143").Value = 3.25
    $ws.Range("O143").Value = 2.1
    $ws.Range("P143").Value = 0.25
    $ws.Range("Q143").Value = 1.95
    $ws.Range("R143").Value = 1.85
    $ws.Range("S143").Value = 2.25
    $ws.Range("T143").Value = 1.775
    $ws.Range("U143").Value = 2.025
    $ws.Range("V143").Value = -1
    $ws.Range("W143").Value = -1
    $ws.Range("X143").Value = 1.1
    $ws.Range("Y143").Value = -1
    $ws.Range("Z143").Value = 0.8500000000000001
    $ws.Range("AA143").Value = -1
    $ws.Range("AB143").Value = 1.025

    # Row 145
    $ws.Range("B145").Value = 7493433
    $ws.Range("E145").Value = 'Sportivo Luqueno'
    $ws.Range("F145").Value = 'Nacional Asuncion'
    $ws.Range("G145").Value = 1
    $ws.Range("H145").Value = 1
    $ws.Range("I145").Value = 'D'
    $ws.Range("J145").Value = 2.75
    $ws.Range("K145").Value = 3.2
    $ws.Range("L145").Value = 2.4
    $ws.Range("M145").Value = 2.75
    $ws.Range("N145").Value = 3.1
    $ws.Range("O145").Value = 2.45
    $ws.Range("P145").Value = 0.25
    $ws.Range("Q145").Value = 1.75
    $ws.Range("R145").Value = 2.05
    $ws.Range("S145").Value = 2.25
    $ws.Range("T145").Value = 2
    $ws.Range("U145").Value = 1.8
    $ws.Range("V145").Value = -1
    $ws.Range("W145").Value = 2.1
    $ws.Range("X145").Value = -1
    $ws.Range("Y145").Value = 0.375
    $ws.Range("Z145").Value = -0.5
    $ws.Range("AA145").Value = -0.5
    $ws.Range("AB145").Value = 0.4


    # --- Rows 231,232: same-date fixtures swapped by source scraper ---

    # Row 231
    $ws.Range("B231").Value = 7609161
    $ws.Range("E231").Value = 'Guarani Asuncion'
    $ws.Range("F231").Value = 'Nacional Asuncion'
    $ws.Range("G231").Value = 3
    $ws.Range("H231").Value = 1
    $ws.Range("I231").Value = 'H'
    $ws.Range("J231").Value = 2.1
    $ws.Range("K231").Value = 3.25
    $ws.Range("L231").Value = 3.6
    $ws.Range("M231").Value = 2.25
    $ws.Range("N231").Value = 3.1
    $ws.Range("O231").Value = 3.4
    $ws.Range("P231").Value = -0.25
    $ws.Range("Q231").Value = 1.9
    $ws.Range("R231").Value = 1.9
    $ws.Range("S231").Value = 2.25
    $ws.Range("T231").Value = 2
    $ws.Range("U231").Value = 1.8
    $ws.Range("V231").Value = 1.25
    $ws.Range("W231").Value = -1
    $ws.Range("X231").Value = -1
    $ws.Range("Y231").Value = 0.8999999999999999
    $ws.Range("Z231").Value = -1
    $ws.Range("AA231").Value = 1
    $ws.Range("AB231").Value = -1

    # Row 232
    $ws.Range("B232").Value = 7609668
    $ws.Range("E232").Value = '2 de Mayo'
    $ws.Range("F232").Value = 'Libertad Asuncion'
    $ws.Range("G232").Value = 2
    $ws.Range("H232").Value = 0
    $ws.Range("I232").Value = 'H'
    $ws.Range("J232").Value = 4.2
    $ws.Range("K232").Value = 3.5
    $ws.Range("L232").Value = 1.85
    $ws.Range("M232").Value = 4
    $ws.Range("N232").Value = 3.4
    $ws.Range("O232").Value = 1.909
    $ws.Range("P232").Value = 0.5
    $ws.Range("Q232").Value = 1.9
    $ws.Range("R232").Value = 1.9
    $ws.Range("S232").Value = 2.25
    $ws.Range("T232").Value = 1.85
    $ws.Range("U232").Value = 1.95
    $ws.Range("V232").Value = 3
    $ws.Range("W232").Value = -1
    $ws.Range("X232").Value = -1
    $ws.Range("Y232").Value = 0.8999999999999999
    $ws.Range("Z232").Value = -1
    $ws.Range("AA232").Value = -0.5
    $ws.Range("AB232").Value = 0.475


    # --- Rows 247-249: live odds refresh for upcoming fixtures ---

    # Row 247
    $ws.Range("M247").Value = 2.4
    $ws.Range("N247").Value = 3.25
    $ws.Range("O247").Value = 2.8
    $ws.Range("Q247").Value = 1.775
    $ws.Range("R247").Value = 2.025
    $ws.Range("S247").Value = 2.25
    $ws.Range("T247").Value = 1.85
    $ws.Range("U247").Value = 1.95

    # Row 248
    $ws.Range("M248").Value = 2.55
    $ws.Range("O248").Value = 2.7
    $ws.Range("Q248").Value = 1.8
    $ws.Range("R248").Value = 2
    $ws.Range("T248").Value = 2
    $ws.Range("U248").Value = 1.8

    # Row 249
    $ws.Range("M249").Value = 2.375
    $ws.Range("N249").Value = 3.25
    $ws.Range("O249").Value = 2.9
    $ws.Range("Q249").Value = 2.05
    $ws.Range("R249").Value = 1.75
    $ws.Range("T249").Value = 1.825
    $ws.Range("U249").Value = 1.975
